$d = $word.ActiveDocument

$rng = $d.Paragraphs(6).Range
$rng.Find.Execute('Munir o aluno de ferramentas para: especificar materiais para aplicação na indústria química e definir os meios possíveis de processar os materiais comuns a sua área de atuação.', $true, $false, $false, $false, $false, $true, 1, $false, 'Introdução aos materiais para a indústria química, propriedades, especificações, seleção, fabricação, aplicação e corrosão.', 2) | Out-Null
$rng = $d.Paragraphs(8).Range
$rng.Find.Execute('144651 - Antonio Fernando Sartori', $true, $false, $false, $false, $false, $true, 1, $false, 'Munir o aluno de ferramentas para: especificar materiais para aplicação na indústria química e definir os meios possíveis de processar os materiais comuns a sua área de atuação.', 2) | Out-Null
$rng = $d.Paragraphs(8).Range
$rng.Find.Execute('3586455 - Cassius Olivio Figueiredo Terra Ruchert', $true, $false, $false, $false, $false, $true, 1, $false, 'Introdução aos materiais. - Seleção de materiais. - Fatores que influenciam na seleção dos materiais (indústria química, petroquímica, Nuclear e outras), melhoria das propriedades mecânicas dos materiais. - Falhas em serviço e em processo.  Produtos siderúrgicos para aplicação em indústrias químicas - Aços carbono e especiais - Ferro fundido. – Processo de fabricação de aços e ferros fundidos, especificações, propriedades e aplicações. Metais e ligas não ferrosas e não metálicas: especificações, propriedades e aplicações. Introdução à corrosão aplicada a engenharia. Pilha Eletroquímica e eletrolítica, meios corrosivos, causas e formas de corrosão, corrosão seletiva, induzida por micromecanismos (MIQ), puntiforme, filiforme, frestas, CST etc... Proteção de superfícies metálicas contra a corrosão, tipo de revestimentos como aspersão térmica, PVD, QVD, etc..', 2) | Out-Null
$rng = $d.Paragraphs(8).Range
$rng.Find.Execute('5840963 - Daniela Camargo Vernilli', $true, $false, $false, $false, $false, $true, 1, $false, 'De acordo com a atual ementa da disciplina propõe-se o uso de uma nova metodologia de ensino com o intuito de abordar o conteúdo de forma mais prática e contextualizada para que o aluno consiga relacionar os conhecimentos teóricos vistos em sala de aula com as outras disciplinas do curso. Assim, avaliação do aluno será feita através de uma prova escrita e por uma apresentação final com base nas atividades práticas desenvolvidas.', 2) | Out-Null
$rng = $d.Paragraphs(8).Range
$rng.Find.Execute('984972 - Hugo Ricardo Zschommler Sandim', $true, $false, $false, $false, $false, $true, 1, $false, 'A nota final será calculada como descrita a seguir: NF= (0,4*Avaliação escrita + 0,6 *Apresentação final)', 2) | Out-Null
$rng = $d.Paragraphs(10).Range
$rng.Find.Execute('Introdução aos materiais para a indústria química, propriedades, especificações, seleção, fabricação, aplicação e corrosão.', $true, $false, $false, $false, $false, $true, 1, $false, 'Devido a cunho prático da disciplina não haverá recuperação.', 2) | Out-Null
$rng = $d.Paragraphs(12).Range
$rng.Find.Execute('Introdução aos materiais. - Seleção de materiais. - Fatores que influenciam na seleção dos materiais (indústria química, petroquímica, Nuclear e outras), melhoria das propriedades mecânicas dos materiais. - Falhas em serviço e em processo.  Produtos siderúrgicos para aplicação em indústrias químicas - Aços carbono e especiais - Ferro fundido. – Processo de fabricação de aços e ferros fundidos, especificações, propriedades e aplicações. Metais e ligas não ferrosas e não metálicas: especificações, propriedades e aplicações. Introdução à corrosão aplicada a engenharia. Pilha Eletroquímica e eletrolítica, meios corrosivos, causas e formas de corrosão, corrosão seletiva, induzida por micromecanismos (MIQ), puntiforme, filiforme, frestas, CST etc... Proteção de superfícies metálicas contra a corrosão, tipo de revestimentos como aspersão térmica, PVD, QVD, etc..', $true, $false, $false, $false, $false, $true, 1, $false, '1)Telles, P. C. S. - Materiais para Equipamentos de Processo - Ed. Interciência, 4º Ed., 1989.^l2)Bresciani, F., E. - Seleção de Materiais Metálicos - Ed. da UNICAMP, 2º Ed.^l3)Freire, J. M. -Materiais de Construção Mecânica - Ed. Livros Técnicos e Científicos, Editora 1993.^l4)A. Remy/ M. Gay/ R. Gonthier - Materiais - Hemus Editora Limitada - 2ª Edição.^l5)Chiaverini, V.Tecnologia Mecânica - Materiais de Construção Mecânica - Vol. II - Ed. McGraw Hill do Brasil Ltda.^l6)Gentil, V. - Corrosão. - Ed. Guanabara Dois, 1982.', 2) | Out-Null
$rng = $d.Paragraphs(14).Range
$rng.Find.Execute('De acordo com a atual ementa da disciplina propõe-se o uso de uma nova metodologia de ensino com o intuito de abordar o conteúdo de forma mais prática e contextualizada para que o aluno consiga relacionar os conhecimentos teóricos vistos em sala de aula com as outras disciplinas do curso. Assim, avaliação do aluno será feita através de uma prova escrita e por uma apresentação final com base nas atividades práticas desenvolvidas.', $true, $false, $false, $false, $false, $true, 1, $false, '144651 - Antonio Fernando Sartori', 2) | Out-Null
$rng = $d.Paragraphs(14).Range
$rng.Find.Execute('A nota final será calculada como descrita a seguir: NF= (0,4*Avaliação escrita + 0,6 *Apresentação final)', $true, $false, $false, $false, $false, $true, 1, $false, '3586455 - Cassius Olivio Figueiredo Terra Ruchert', 2) | Out-Null
$rng = $d.Paragraphs(14).Range
$rng.Find.Execute('Devido a cunho prático da disciplina não haverá recuperação.', $true, $false, $false, $false, $false, $true, 1, $false, '5840963 - Daniela Camargo Vernilli', 2) | Out-Null
$rng = $d.Paragraphs(16).Range
$rng.Find.Execute('1)Telles, P. C. S. - Materiais para Equipamentos de Processo - Ed. Interciência, 4º Ed., 1989.^l2)Bresciani, F., E. - Seleção de Materiais Metálicos - Ed. da UNICAMP, 2º Ed.^l3)Freire, J. M. -Materiais de Construção Mecânica - Ed. Livros Técnicos e Científicos, Editora 1993.^l4)A. Remy/ M. Gay/ R. Gonthier - Materiais - Hemus Editora Limitada - 2ª Edição.^l5)Chiaverini, V.Tecnologia Mecânica - Materiais de Construção Mecânica - Vol. II - Ed. McGraw Hill do Brasil Ltda.^l6)Gentil, V. - Corrosão. - Ed. Guanabara Dois, 1982.', $true, $false, $false, $false, $false, $true, 1, $false, '984972 - Hugo Ricardo Zschommler Sandim', 2) | Out-Null
